$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 119
$ws.Range("A119").Value = '2024-10-28 19:27:52'
$ws.Range("B119").Value = 'ORM-0516056'
$ws.Range("C119").Value = 'MCULE-2227031507'
$ws.Range("D119").Value = 'C21H17CLN2O2'
$ws.Range("E119").Value = 'Success'
$ws.Range("F119").Value = '/home/robekott/ERAT/MCULE/mcule_test.sdf'
$ws.Range("G119").Value = 201
$ws.Range("H119").Value = '{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671fc997b0480736afa1ed64"},"data":{"type":"material","id":"asset:671fc997b0480736afa1ed64","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671fc997b0480736afa1ed64"},"attributes":{"library":"Compounds","assetTypeId":"5d6e0287ee35880008c18db5","assetId":"671fc997b0480736afa1ed64","id":"asset:671fc997b0480736afa1ed64","eid":"asset:671fc997b0480736afa1ed64","name":"ORM-0516056","synonyms":["C(=O)(NNC(=O)CC1CCC(C2CCCCC2)CC1)C1CCC(CC1)CL","C21H17CLN2O2"],"description":"","createdAt":"2024-10-28T17:27:51.830Z","editedAt":"2024-10-28T17:27:51.830Z","type":"asset","digest":"66868543","fields":{"Chemical Name":{"value":"N''-(2-([1,1''-biphenyl]-4-yl)acetyl)-4-chlorobenzohydrazide"},"Description":{"value":""},"Exact Mass":{"value":"364.09786"},"Material Library Type":{"value":"Compounds"},"Molecular Formula":{"value":"C<sub>21</sub>H<sub>17</sub>ClN<sub>2</sub>O<sub>2</sub>"},"Molecular Weight":{"value":"364.83 g/mol"},"Name":{"value":"ORM-0516056"},"Stereochemistry":{"value":"No stereochemistry"}},"flags":{"canTrash":true}},"relationships":{"batches":{"data":[{"type":"material","id":"batch:671fc998b0480736afa1ed65","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:671fc998b0480736afa1ed65"}}}]},"ancestors":{"data":[{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"}}}]},"chemicalDrawing":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671fc997b0480736afa1ed64/drawing?format=cdxml"},"data":{"type":"materialDrawing","id":"asset:671fc997b0480736afa1ed64"}},"createdBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"editedBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"owner":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}}}},"included":[{"type":"material","id":"batch:671fc998b0480736afa1ed65","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:671fc998b0480736afa1ed65"},"attributes":{"type":"batch","eid":"batch:671fc998b0480736afa1ed65","name":"ORM-0516056-001","digest":"83571980","fields":{"Batch Chemical Name":{"value":"N''-(2-([1,1''-biphenyl]-4-yl)acetyl)-4-chlorobenzohydrazide"},"Batch Molecular Formula":{"value":"C<sub>21</sub>H<sub>17</sub>ClN<sub>2</sub>O<sub>2</sub>"},"Batch Purpose":{"value":"Dummy compound"},"Batch Type":{"value":"Discovery"},"Chemist":{"value":"TestUser MCChemist"},"Description":{"value":""},"Formula Mass":{"value":"364.83 g/mol"},"Library ID":{"value":""},"Name":{"value":"ORM-0516056-001"},"Project":{"value":"Unspecified"},"Source":{"value":"Internal"},"Submission Date":{"value":"2024-10-28T17:27:52.147Z"},"Submitter":{"value":"{userId=114, userName=timo.kangasperko@orion.fi, flags={isSystemStandardUser=true}, alias=timoka, email=timo.kangasperko@orion.fi, firstName=Timo, lastName=Kangasperko, picture={}, isEnabled=true}"},"Supplier Name":{"value":"Chembridge"},"Synthesis Date":{"value":"2011-10-10T14:48Z"}}}},{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"},"attributes":{"type":"assetType","eid":"assetType:5d6e0287ee35880008c18db5","name":"Compounds","digest":"59951413","fields":{}}},{"type":"user","id":"114","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"attributes":{"userId":"114","userName":"timo.kangasperko@orion.fi","flags":{"isSystemStandardUser":true},"alias":"timoka","email":"timo.kangasperko@orion.fi","firstName":"Timo","lastName":"Kangasperko","isEnabled":true},"relationships":{"systemGroups":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114/systemGroups"}}}},{"type":"materialDrawing","id":"asset:671fc997b0480736afa1ed64","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671fc997b0480736afa1ed64/drawing?format=cdxml"},"attributes":{"id":"asset:671fc997b0480736afa1ed64","type":"CHEMICAL_DRAWING"}}]}'

# Row 120
$ws.Range("A120").Value = '2024-10-28 19:27:56'
$ws.Range("B120").Value = 'ORM-0516057'
$ws.Range("C120").Value = 'MCULE-3988458386'
$ws.Range("D120").Value = 'C16H14CLN5O'
$ws.Range("E120").Value = 'Success'
$ws.Range("F120").Value = '/home/robekott/ERAT/MCULE/mcule_test.sdf'
$ws.Range("G120").Value = 201
$ws.Range("H120").Value = '{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671fc99b1a12207a92bcb05a"},"data":{"type":"material","id":"asset:671fc99b1a12207a92bcb05a","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671fc99b1a12207a92bcb05a"},"attributes":{"library":"Compounds","assetTypeId":"5d6e0287ee35880008c18db5","assetId":"671fc99b1a12207a92bcb05a","id":"asset:671fc99b1a12207a92bcb05a","eid":"asset:671fc99b1a12207a92bcb05a","name":"ORM-0516057","synonyms":["C1(C([NH]NN1)NC1CCC(CL)CC1)C(=O)NC1CCC(CC1)C","C16H14CLN5O"],"description":"","createdAt":"2024-10-28T17:27:55.840Z","editedAt":"2024-10-28T17:27:55.840Z","type":"asset","digest":"93385593","fields":{"Chemical Name":{"value":"5-((4-chlorophenyl)amino)-N-(p-tolyl)-1H-1,2,3-triazole-4-carboxamide"},"Description":{"value":""},"Exact Mass":{"value":"327.08869"},"Material Library Type":{"value":"Compounds"},"Molecular Formula":{"value":"C<sub>16</sub>H<sub>14</sub>ClN<sub>5</sub>O"},"Molecular Weight":{"value":"327.77 g/mol"},"Name":{"value":"ORM-0516057"},"Stereochemistry":{"value":"No stereochemistry"}},"flags":{"canTrash":true}},"relationships":{"batches":{"data":[{"type":"material","id":"batch:671fc99c1a12207a92bcb05b","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:671fc99c1a12207a92bcb05b"}}}]},"ancestors":{"data":[{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"}}}]},"chemicalDrawing":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671fc99b1a12207a92bcb05a/drawing?format=cdxml"},"data":{"type":"materialDrawing","id":"asset:671fc99b1a12207a92bcb05a"}},"createdBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"editedBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"owner":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}}}},"included":[{"type":"material","id":"batch:671fc99c1a12207a92bcb05b","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:671fc99c1a12207a92bcb05b"},"attributes":{"type":"batch","eid":"batch:671fc99c1a12207a92bcb05b","name":"ORM-0516057-001","digest":"59455064","fields":{"Batch Chemical Name":{"value":"5-((4-chlorophenyl)amino)-N-(p-tolyl)-1H-1,2,3-triazole-4-carboxamide"},"Batch Molecular Formula":{"value":"C<sub>16</sub>H<sub>14</sub>ClN<sub>5</sub>O"},"Batch Purpose":{"value":"Dummy compound"},"Batch Type":{"value":"Discovery"},"Chemist":{"value":"TestUser MCChemist"},"Description":{"value":""},"Formula Mass":{"value":"327.77 g/mol"},"Library ID":{"value":""},"Name":{"value":"ORM-0516057-001"},"Project":{"value":"Unspecified"},"Source":{"value":"Internal"},"Submission Date":{"value":"2024-10-28T17:27:56.113Z"},"Submitter":{"value":"{userId=114, userName=timo.kangasperko@orion.fi, flags={isSystemStandardUser=true}, alias=timoka, email=timo.kangasperko@orion.fi, firstName=Timo, lastName=Kangasperko, picture={}, isEnabled=true}"},"Supplier Name":{"value":"ChemDiv"},"Synthesis Date":{"value":"2011-10-10T14:48Z"}}}},{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"},"attributes":{"type":"assetType","eid":"assetType:5d6e0287ee35880008c18db5","name":"Compounds","digest":"59951413","fields":{}}},{"type":"materialDrawing","id":"asset:671fc99b1a12207a92bcb05a","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671fc99b1a12207a92bcb05a/drawing?format=cdxml"},"attributes":{"id":"asset:671fc99b1a12207a92bcb05a","type":"CHEMICAL_DRAWING"}},{"type":"user","id":"114","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"attributes":{"userId":"114","userName":"timo.kangasperko@orion.fi","flags":{"isSystemStandardUser":true},"alias":"timoka","email":"timo.kangasperko@orion.fi","firstName":"Timo","lastName":"Kangasperko","isEnabled":true},"relationships":{"systemGroups":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114/systemGroups"}}}}]}'

Write-Output "done"